# feat: add 2022-Q1 data
#
# - Insert a new sheet "2022-Q1" (with per-fund holdings) between the
#   existing "2021-Q4" and "总计" sheets.
# - Insert a new summary row for "2022-Q1" at the top of the "总计"
#   (totals) sheet's data, pushing the 2021-Q4 / 2021-Q3 rows down.

$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet right after "2021-Q4".
#    NOTE: this shifts the index of every sheet after it (in particular
#    "总计"), so any sheet handle obtained before this point other than
#    "2021-Q3"/"2021-Q4" (which sit before the insertion point) must be
#    re-fetched afterwards instead of reused.
# ---------------------------------------------------------------------
$new = $wb.Worksheets.Add($null, $q4)
$new.Name = "2022-Q1"

# Match the page-margin metadata used by the other per-quarter sheets.
$new.PageSetup.LeftMargin = $q4.PageSetup.LeftMargin
$new.PageSetup.RightMargin = $q4.PageSetup.RightMargin
$new.PageSetup.TopMargin = $q4.PageSetup.TopMargin
$new.PageSetup.BottomMargin = $q4.PageSetup.BottomMargin
$new.PageSetup.HeaderMargin = $q4.PageSetup.HeaderMargin
$new.PageSetup.FooterMargin = $q4.PageSetup.FooterMargin

# Copy the header-row / index-column formatting from "2021-Q4" so the new
# sheet's look matches its siblings (bold header with border, s=2 style).
$q4.Range("B1:H1").Copy()
$new.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2").Copy()
$new.Range("A2:A3").PasteSpecial(-4122)

# Header row.
$new.Range("B1").Value = "基金代码"
$new.Range("C1").Value = "基金名称"
$new.Range("D1").Value = "基金规模"
$new.Range("E1").Value = "股票总仓位"
$new.Range("F1").Value = "仓位占比"
$new.Range("G1").Value = "持有市值(亿元)"
$new.Range("H1").Value = "仓位排名"

# Columns B (fund code, to preserve the leading zero) and D/E/F/G hold
# text-formatted numbers (mirrors how the other quarter sheets store
# them) while H is a real number. Force "@" (text) number format before
# assignment so the value lands as a string instead of Excel's
# auto-numeric coercion, then drop the format again so no stray style
# survives on the cell.
$new.Range("A2").Value = 0
$new.Range("B2").NumberFormat = "@"
$new.Range("B2").Value = "009753"
$new.Range("B2").ClearFormats()
$new.Range("C2").Value = "中欧美益稳健两年持有期混合A"
$new.Range("D2:G2").NumberFormat = "@"
$new.Range("D2").Value = "2.63"
$new.Range("E2").Value = "23.02"
$new.Range("F2").Value = "1.39"
$new.Range("G2").Value = "0.0366"
$new.Range("D2:G2").ClearFormats()
$new.Range("H2").Value = 3

$new.Range("A3").Value = 1
$new.Range("B3").NumberFormat = "@"
$new.Range("B3").Value = "009754"
$new.Range("B3").ClearFormats()
$new.Range("C3").Value = "中欧美益稳健两年持有期混合C"
$new.Range("D3:G3").NumberFormat = "@"
$new.Range("D3").Value = "0.23"
$new.Range("E3").Value = "23.02"
$new.Range("F3").Value = "1.39"
$new.Range("G3").Value = "0.0032"
$new.Range("D3:G3").ClearFormats()
$new.Range("H3").Value = 3

# ---------------------------------------------------------------------
# 2. Insert the "2022-Q1" summary row at the top of the "总计" sheet.
#    Re-fetch the sheet by name now that "2022-Q1" has shifted its index.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.04

# Renumber the index column for the rows that got pushed down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2

# ---------------------------------------------------------------------
# Restore the originally active sheet/selection.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q3").Activate()
